$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "split PHSMs periods into two stage" - the Dickey-Fuller test / p-value
# results were recomputed for the new two-stage PHSMs periods. Disease
# names and row order are unchanged; only the numeric-looking text values
# in columns B (Dickey_Fuller.Dickey-Fuller) and C (p_value) change.
$data = @(
  @("HAV", "-3.17", "0.0961"),
  @("HEV", "-2.05", "0.5548"),
  @("Typhoid and paratyphoid fever", "-3.09", "0.1212"),
  @("AHC", "-4.43", "0.01"),
  @("Dysentery", "-5.33", "0.01"),
  @("HFMD", "-3.71", "0.0262"),
  @("Other infectious diarrhea", "-3.33", "0.0689"),
  @("Brucellosis", "-2.8", "0.2456"),
  @("Dengue fever", "-4.3", "0.01"),
  @("HFRS", "-3.13", "0.1058"),
  @("Japanese encephalitis", "-4.31", "0.01"),
  @("Malaria", "-4.62", "0.01"),
  @("Hydatidosis", "-3.03", "0.146"),
  @("Typhus", "-4.08", "0.01"),
  @("AIDS", "-3.27", "0.079"),
  @("Gonorrhea", "-0.09", "0.99"),
  @("HBV", "-2", "0.578"),
  @("HCV", "-1.69", "0.7048"),
  @("Syphilis", "-2.38", "0.4177"),
  @("Pertussis", "-0.9", "0.9508"),
  @("Scarlet fever", "-1.87", "0.6293"),
  @("Tuberculosis", "-1.88", "0.6248"),
  @("Mumps", "-1.87", "0.6318"),
  @("Rubella", "-3.49", "0.0458")
)

# Only the Dickey-Fuller (col B) and p_value (col C) figures move; the
# disease names (col A) stay put. Re-stamp every B cell (all 24 change),
# and only touch a C cell when its text actually changes, so cells whose
# p_value text is unchanged (e.g. still "0.01") keep their original
# formatting/style untouched.
$row = 2
foreach ($item in $data) {
    $newB = $item[1]
    $newC = $item[2]
    $oldC = $ws.Cells.Item($row, 3).Text

    # Values are stored as text (not real numbers) in this workbook, so
    # force text formatting before writing numeric-looking strings,
    # otherwise Excel would auto-convert them to numbers.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $newB

    if ("$oldC" -ne $newC) {
        $ws.Cells.Item($row, 3).NumberFormat = "@"
        $ws.Cells.Item($row, 3).Value = $newC
    }

    $row++
}
